# Apply the "Add data for 2022-06-15" update to the carjacking-by-neighborhood-by-month workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the sheet tab (Through 2022-06-06 -> Through 2022-06-07) ---
$ws.Name = "Through 2022-06-07"

# --- 2. Update the running-month header label in B1 ---
$ws.Range("B1").Value = "June 2022 (through June 07)"

# --- 3. Small single-cell increments scattered through the sheet ---
$ws.Range("AF5").Value  = 1
$ws.Range("H10").Value  = 1
$ws.Range("AL10").Value = 2
$ws.Range("B13").Value  = 1
$ws.Range("N14").Value  = 1
$ws.Range("Z16").Value  = 1
$ws.Range("N19").Value  = 1
$ws.Range("AL20").Value = 1
$ws.Range("H59").Value  = 1
$ws.Range("B68").Value  = 1
$ws.Range("AL85").Value = 1
$ws.Range("N94").Value  = 2
$ws.Range("B96").Value  = 1
$ws.Range("Z96").Value  = 1

# --- 4. Rows 24-27: neighborhood rows shifted because three new rows were
#        inserted ahead of "United Center" in the underlying neighborhood
#        list, and "United Center" itself received updated counts.
#        Row 24 = Chicago Lawn, Row 25 = Ashburn, Row 26 = Rogers Park,
#        Row 27 = United Center (with revised counts).

$ws.Range("A24").Value = "Chicago Lawn"
$ws.Range("A25").Value = "Ashburn"
$ws.Range("A26").Value = "Rogers Park"
$ws.Range("A27").Value = "United Center"

$row24 = New-Object 'object[,]' 1,48
$row24[0,0]  = ""
$row24[0,1]  = 2
$row24[0,2]  = 5
$row24[0,3]  = 4
$row24[0,4]  = 5
$row24[0,5]  = 3
$row24[0,6]  = ""
$row24[0,7]  = ""
$row24[0,8]  = ""
$row24[0,9]  = 1
$row24[0,10] = ""
$row24[0,11] = 6
$row24[0,12] = 2
$row24[0,13] = 1
$row24[0,14] = ""
$row24[0,15] = 1
$row24[0,16] = 1
$row24[0,17] = 3
$row24[0,18] = ""
$row24[0,19] = ""
$row24[0,20] = 2
$row24[0,21] = 1
$row24[0,22] = 2
$row24[0,23] = 1
$row24[0,24] = ""
$row24[0,25] = 3
$row24[0,26] = ""
$row24[0,27] = 3
$row24[0,28] = 2
$row24[0,29] = 4
$row24[0,30] = ""
$row24[0,31] = 2
$row24[0,32] = ""
$row24[0,33] = 1
$row24[0,34] = 2
$row24[0,35] = ""
$row24[0,36] = ""
$row24[0,37] = 2
$row24[0,38] = ""
$row24[0,39] = 2
$row24[0,40] = 4
$row24[0,41] = ""
$row24[0,42] = ""
$row24[0,43] = ""
$row24[0,44] = ""
$row24[0,45] = 1
$row24[0,46] = ""
$row24[0,47] = ""
$ws.Range("B24:AW24").Value = $row24

$row25 = New-Object 'object[,]' 1,48
$row25[0,0]  = ""
$row25[0,1]  = 2
$row25[0,2]  = ""
$row25[0,3]  = ""
$row25[0,4]  = ""
$row25[0,5]  = 1
$row25[0,6]  = ""
$row25[0,7]  = ""
$row25[0,8]  = 3
$row25[0,9]  = ""
$row25[0,10] = 1
$row25[0,11] = ""
$row25[0,12] = ""
$row25[0,13] = ""
$row25[0,14] = ""
$row25[0,15] = ""
$row25[0,16] = ""
$row25[0,17] = ""
$row25[0,18] = ""
$row25[0,19] = 1
$row25[0,20] = ""
$row25[0,21] = ""
$row25[0,22] = ""
$row25[0,23] = ""
$row25[0,24] = ""
$row25[0,25] = ""
$row25[0,26] = 1
$row25[0,27] = 1
$row25[0,28] = 2
$row25[0,29] = 3
$row25[0,30] = ""
$row25[0,31] = 1
$row25[0,32] = 2
$row25[0,33] = ""
$row25[0,34] = 1
$row25[0,35] = ""
$row25[0,36] = ""
$row25[0,37] = ""
$row25[0,38] = ""
$row25[0,39] = ""
$row25[0,40] = ""
$row25[0,41] = 1
$row25[0,42] = ""
$row25[0,43] = ""
$row25[0,44] = 1
$row25[0,45] = ""
$row25[0,46] = ""
$row25[0,47] = ""
$ws.Range("B25:AW25").Value = $row25

$row26 = New-Object 'object[,]' 1,48
$row26[0,0]  = ""
$row26[0,1]  = 1
$row26[0,2]  = ""
$row26[0,3]  = 2
$row26[0,4]  = 6
$row26[0,5]  = 1
$row26[0,6]  = ""
$row26[0,7]  = 1
$row26[0,8]  = 1
$row26[0,9]  = 1
$row26[0,10] = ""
$row26[0,11] = 2
$row26[0,12] = ""
$row26[0,13] = ""
$row26[0,14] = 1
$row26[0,15] = ""
$row26[0,16] = ""
$row26[0,17] = ""
$row26[0,18] = ""
$row26[0,19] = ""
$row26[0,20] = ""
$row26[0,21] = ""
$row26[0,22] = ""
$row26[0,23] = 1
$row26[0,24] = ""
$row26[0,25] = ""
$row26[0,26] = 1
$row26[0,27] = ""
$row26[0,28] = ""
$row26[0,29] = 1
$row26[0,30] = ""
$row26[0,31] = ""
$row26[0,32] = ""
$row26[0,33] = ""
$row26[0,34] = ""
$row26[0,35] = ""
$row26[0,36] = ""
$row26[0,37] = ""
$row26[0,38] = 1
$row26[0,39] = 2
$row26[0,40] = ""
$row26[0,41] = ""
$row26[0,42] = ""
$row26[0,43] = ""
$row26[0,44] = ""
$row26[0,45] = ""
$row26[0,46] = ""
$row26[0,47] = 1
$ws.Range("B26:AW26").Value = $row26

$row27 = New-Object 'object[,]' 1,48
$row27[0,0]  = ""
$row27[0,1]  = 1
$row27[0,2]  = 1
$row27[0,3]  = ""
$row27[0,4]  = 2
$row27[0,5]  = 5
$row27[0,6]  = ""
$row27[0,7]  = 1
$row27[0,8]  = ""
$row27[0,9]  = ""
$row27[0,10] = 5
$row27[0,11] = 4
$row27[0,12] = ""
$row27[0,13] = 2
$row27[0,14] = ""
$row27[0,15] = ""
$row27[0,16] = 1
$row27[0,17] = ""
$row27[0,18] = ""
$row27[0,19] = 1
$row27[0,20] = ""
$row27[0,21] = 1
$row27[0,22] = ""
$row27[0,23] = 3
$row27[0,24] = ""
$row27[0,25] = 1
$row27[0,26] = ""
$row27[0,27] = 1
$row27[0,28] = 1
$row27[0,29] = 2
$row27[0,30] = 1
$row27[0,31] = 2
$row27[0,32] = ""
$row27[0,33] = ""
$row27[0,34] = ""
$row27[0,35] = 2
$row27[0,36] = ""
$row27[0,37] = 1
$row27[0,38] = 2
$row27[0,39] = 1
$row27[0,40] = ""
$row27[0,41] = 2
$row27[0,42] = ""
$row27[0,43] = ""
$row27[0,44] = ""
$row27[0,45] = ""
$row27[0,46] = ""
$row27[0,47] = ""
$ws.Range("B27:AW27").Value = $row27
